$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Experimental -> "false"
# Using a quoted-text formula + paste-values round-trip forces the
# literal text "false" to be stored as a shared string instead of
# being auto-coerced into an Excel boolean TRUE/FALSE cell.
$cell = $ws.Range("B7")
$cell.Formula = "=""false"""
$cell.Copy()
$cell.PasteSpecial(-4163, 0)

# Date -> updated generation timestamp
$ws.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# Description -> populated (falls back to the Title text)
$ws.Range("B12").Value = "Assertion of Phenotypic Feature Codes"
